$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "codeforiati:group-name" (col C) and "codeforiati:group-code" (col D)
# columns - including their header cells - swap places, so that the code
# (e.g. "BE") now appears in column C and the name (e.g. "Belgique (la)")
# appears in column D.
$used = $ws.UsedRange
$lastRow = $used.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cCell = $ws.Cells.Item($r, 3)
    $dCell = $ws.Cells.Item($r, 4)
    $cVal = $cCell.Value2
    $dVal = $dCell.Value2
    $cCell.Value2 = $dVal
    $dCell.Value2 = $cVal
}
